$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrainingRun2")

# Row 22
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = 8.037410736083984
$ws.Cells.Item(22, 4).Value = 0.9009041786193848
$ws.Cells.Item(22, 5).Value = 0.9038553721632818
$ws.Cells.Item(22, 6).Value = 0.8642713813112477
$ws.Cells.Item(22, 7).Value = 0.8632959663119331

# Row 23
$ws.Cells.Item(23, 1).Value = 2
$ws.Cells.Item(23, 2).Value = 17.80628418922424
$ws.Cells.Item(23, 3).Value = 9.768873453140259
$ws.Cells.Item(23, 4).Value = 0.9342758655548096
$ws.Cells.Item(23, 5).Value = 0.9317103189859292
$ws.Cells.Item(23, 6).Value = 0.8955077709887093
$ws.Cells.Item(23, 7).Value = 0.89482506946225

# Row 24
$ws.Cells.Item(24, 1).Value = 3
$ws.Cells.Item(24, 2).Value = 26.69517970085144
$ws.Cells.Item(24, 3).Value = 8.888895511627197
$ws.Cells.Item(24, 4).Value = 0.9247410893440247
$ws.Cells.Item(24, 5).Value = 0.9290847545054026
$ws.Cells.Item(24, 6).Value = 0.8876495648344066
$ws.Cells.Item(24, 7).Value = 0.8862864352204776

# Row 25
$ws.Cells.Item(25, 1).Value = 4
$ws.Cells.Item(25, 2).Value = 35.97301006317139
$ws.Cells.Item(25, 3).Value = 9.277830362319946
$ws.Cells.Item(25, 4).Value = 0.959855318069458
$ws.Cells.Item(25, 5).Value = 0.9572491923417846
$ws.Cells.Item(25, 6).Value = 0.9222882247224139
$ws.Cells.Item(25, 7).Value = 0.9222511348001811

# Row 26
$ws.Cells.Item(26, 1).Value = 5
$ws.Cells.Item(26, 2).Value = 44.3629789352417
$ws.Cells.Item(26, 3).Value = 8.389968872070312
$ws.Cells.Item(26, 4).Value = 0.9638336300849915
$ws.Cells.Item(26, 5).Value = 0.9607068603829316
$ws.Cells.Item(26, 6).Value = 0.9259423637818004
$ws.Cells.Item(26, 7).Value = 0.9253000781939975

# Row 27
$ws.Cells.Item(27, 1).Value = 6
$ws.Cells.Item(27, 2).Value = 52.88276696205139
$ws.Cells.Item(27, 3).Value = 8.519788026809692
$ws.Cells.Item(27, 4).Value = 0.9624527096748352
$ws.Cells.Item(27, 5).Value = 0.9631393028274935
$ws.Cells.Item(27, 6).Value = 0.9259345001925611
$ws.Cells.Item(27, 7).Value = 0.9267034133505136

# Row 28
$ws.Cells.Item(28, 1).Value = 7
$ws.Cells.Item(28, 2).Value = 62.56203389167786
$ws.Cells.Item(28, 3).Value = 9.679266929626465
$ws.Cells.Item(28, 4).Value = 0.9640966653823853
$ws.Cells.Item(28, 5).Value = 0.9616812555574823
$ws.Cells.Item(28, 6).Value = 0.9252630152870445
$ws.Cells.Item(28, 7).Value = 0.9236578548190372

# Row 29
$ws.Cells.Item(29, 1).Value = 8
$ws.Cells.Item(29, 2).Value = 73.32969045639038
$ws.Cells.Item(29, 3).Value = 10.76765656471252
$ws.Cells.Item(29, 4).Value = 0.9319414496421814
$ws.Cells.Item(29, 5).Value = 0.9387156493625052
$ws.Cells.Item(29, 6).Value = 0.898761260631326
$ws.Cells.Item(29, 7).Value = 0.8970461895495352

# Row 30
$ws.Cells.Item(30, 1).Value = 9
$ws.Cells.Item(30, 2).Value = 82.04559350013733
$ws.Cells.Item(30, 3).Value = 8.715903043746948
$ws.Cells.Item(30, 4).Value = 0.9741575121879578
$ws.Cells.Item(30, 5).Value = 0.972389452267386
$ws.Cells.Item(30, 6).Value = 0.9722806524280484
$ws.Cells.Item(30, 7).Value = 0.9722140580229934

# Row 31
$ws.Cells.Item(31, 1).Value = 10
$ws.Cells.Item(31, 2).Value = 90.77180194854736
$ws.Cells.Item(31, 3).Value = 8.726208448410034
$ws.Cells.Item(31, 4).Value = 0.9760973453521729
$ws.Cells.Item(31, 5).Value = 0.9756835050505298
$ws.Cells.Item(31, 6).Value = 0.9732429083582714
$ws.Cells.Item(31, 7).Value = 0.9739777117311195

# Row 32
$ws.Cells.Item(32, 1).Value = 11
$ws.Cells.Item(32, 2).Value = 100.0763611793518
$ws.Cells.Item(32, 3).Value = 9.304559230804443
$ws.Cells.Item(32, 4).Value = 0.9740259647369385
$ws.Cells.Item(32, 5).Value = 0.971767023505224
$ws.Cells.Item(32, 6).Value = 0.9728566953711947
$ws.Cells.Item(32, 7).Value = 0.9718508223790493

# Row 33
$ws.Cells.Item(33, 1).Value = 12
$ws.Cells.Item(33, 2).Value = 108.8250782489777
$ws.Cells.Item(33, 3).Value = 8.748717069625854
$ws.Cells.Item(33, 4).Value = 0.9755712747573853
$ws.Cells.Item(33, 5).Value = 0.97405812988468
$ws.Cells.Item(33, 6).Value = 0.9743299476204527
$ws.Cells.Item(33, 7).Value = 0.9740008877594032

# Row 34
$ws.Cells.Item(34, 1).Value = 13
$ws.Cells.Item(34, 2).Value = 117.0536303520203
$ws.Cells.Item(34, 3).Value = 8.228552103042603
$ws.Cells.Item(34, 4).Value = 0.9760644435882568
$ws.Cells.Item(34, 5).Value = 0.9737320709625401
$ws.Cells.Item(34, 6).Value = 0.9740349246244088
$ws.Cells.Item(34, 7).Value = 0.9732360897635419

# Row 35
$ws.Cells.Item(35, 1).Value = 14
$ws.Cells.Item(35, 2).Value = 126.663524389267
$ws.Cells.Item(35, 3).Value = 9.609894037246704
$ws.Cells.Item(35, 4).Value = 0.9765576124191284
$ws.Cells.Item(35, 5).Value = 0.976099714504825
$ws.Cells.Item(35, 6).Value = 0.9754770032127669
$ws.Cells.Item(35, 7).Value = 0.9756446323668319

# Row 36
$ws.Cells.Item(36, 1).Value = 15
$ws.Cells.Item(36, 2).Value = 135.897049665451
$ws.Cells.Item(36, 3).Value = 9.233525276184082
$ws.Cells.Item(36, 4).Value = 0.971329927444458
$ws.Cells.Item(36, 5).Value = 0.9709969823862185
$ws.Cells.Item(36, 6).Value = 0.9706319419241851
$ws.Cells.Item(36, 7).Value = 0.9704154821911015

# Row 37
$ws.Cells.Item(37, 1).Value = 16
$ws.Cells.Item(37, 2).Value = 144.0664365291595
$ws.Cells.Item(37, 3).Value = 8.169386863708496
$ws.Cells.Item(37, 4).Value = 0.9687325358390808
$ws.Cells.Item(37, 5).Value = 0.969710145939751
$ws.Cells.Item(37, 6).Value = 0.9694352366053616
$ws.Cells.Item(37, 7).Value = 0.9688816832119176

# Row 38
$ws.Cells.Item(38, 1).Value = 17
$ws.Cells.Item(38, 2).Value = 152.2339911460876
$ws.Cells.Item(38, 3).Value = 8.167554616928101
$ws.Cells.Item(38, 4).Value = 0.9661351442337036
$ws.Cells.Item(38, 5).Value = 0.9655260876826288
$ws.Cells.Item(38, 6).Value = 0.9617594163118518
$ws.Cells.Item(38, 7).Value = 0.9601438496307407

# Row 39
$ws.Cells.Item(39, 1).Value = 18
$ws.Cells.Item(39, 2).Value = 160.7514872550964
$ws.Cells.Item(39, 3).Value = 8.517496109008789
$ws.Cells.Item(39, 4).Value = 0.9788591265678406
$ws.Cells.Item(39, 5).Value = 0.9780127278479415
$ws.Cells.Item(39, 6).Value = 0.978218863794879
$ws.Cells.Item(39, 7).Value = 0.9780014459550754

# Row 40
$ws.Cells.Item(40, 1).Value = 19
$ws.Cells.Item(40, 2).Value = 171.4699847698212
$ws.Cells.Item(40, 3).Value = 10.71849751472473
$ws.Cells.Item(40, 4).Value = 0.9827058911323547
$ws.Cells.Item(40, 5).Value = 0.9816476376908092
$ws.Cells.Item(40, 6).Value = 0.980898479912813
$ws.Cells.Item(40, 7).Value = 0.9810578176394261

# Row 41
$ws.Cells.Item(41, 1).Value = 20
$ws.Cells.Item(41, 2).Value = 180.6203641891479
$ws.Cells.Item(41, 3).Value = 9.150379419326782
$ws.Cells.Item(41, 4).Value = 0.9794509410858154
$ws.Cells.Item(41, 5).Value = 0.977603868229227
$ws.Cells.Item(41, 6).Value = 0.9768589411552572
$ws.Cells.Item(41, 7).Value = 0.9762675116853456

# Row 42
$ws.Cells.Item(42, 1).Value = 1
$ws.Cells.Item(42, 2).Value = 20.5004358291626
$ws.Cells.Item(42, 4).Value = 0.9069867134094238
$ws.Cells.Item(42, 5).Value = 0.9110572373172203
$ws.Cells.Item(42, 6).Value = 0.8703144659650774
$ws.Cells.Item(42, 7).Value = 0.8709239189736855

# Row 43
$ws.Cells.Item(43, 1).Value = 2
$ws.Cells.Item(43, 2).Value = 44.49518322944641
$ws.Cells.Item(43, 3).Value = 23.99474740028381
$ws.Cells.Item(43, 4).Value = 0.9361499547958374
$ws.Cells.Item(43, 5).Value = 0.9370535028960167
$ws.Cells.Item(43, 6).Value = 0.8979039293265859
$ws.Cells.Item(43, 7).Value = 0.8990548970545544

# Row 44
$ws.Cells.Item(44, 1).Value = 3
$ws.Cells.Item(44, 2).Value = 67.79964804649353
$ws.Cells.Item(44, 3).Value = 23.30446481704712
$ws.Cells.Item(44, 4).Value = 0.9507151246070862
$ws.Cells.Item(44, 5).Value = 0.9513855397267067
$ws.Cells.Item(44, 6).Value = 0.9118300609374642
$ws.Cells.Item(44, 7).Value = 0.9128254086923384

# Row 45
$ws.Cells.Item(45, 1).Value = 4
$ws.Cells.Item(45, 2).Value = 92.18682742118835
$ws.Cells.Item(45, 3).Value = 24.38717937469482
$ws.Cells.Item(45, 4).Value = 0.9635705947875977
$ws.Cells.Item(45, 5).Value = 0.9625932551729826
$ws.Cells.Item(45, 6).Value = 0.9270077892056477
$ws.Cells.Item(45, 7).Value = 0.9273513163863312
